$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 2011
$ws.Range("I3").Value = 2128
$ws.Range("G4").Value = 1427
$ws.Range("I4").Value = 531
$ws.Range("I5").Value = 187
$ws.Range("I6").Value = 2533
$ws.Range("G7").Value = 24646
$ws.Range("I7").Value = 7390

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 29
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 25

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 47
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 66
$ws.Range("I3").Value = 98
$ws.Range("I6").Value = 98
$ws.Range("I7").Value = 283

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 21
$ws.Range("I3").Value = 21
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 20
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 48
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I6").Value = 45
$ws.Range("I8").Value = 469
$ws.Range("I10").Value = 58
$ws.Range("I16").Value = 25
$ws.Range("G20").Value = 620
$ws.Range("I23").Value = 62
$ws.Range("I26").Value = 8
$ws.Range("I27").Value = 69
$ws.Range("I29").Value = 477
$ws.Range("I30").Value = 25
$ws.Range("I31").Value = 75
$ws.Range("I33").Value = 349
$ws.Range("I42").Value = 247
$ws.Range("I43").Value = 68
$ws.Range("I45").Value = 14
$ws.Range("I51").Value = 69
$ws.Range("I54").Value = 167
$ws.Range("I60").Value = 40
$ws.Range("I64").Value = 76
$ws.Range("I65").Value = 174
$ws.Range("I67").Value = 283
$ws.Range("I69").Value = 19
$ws.Range("I72").Value = 25
$ws.Range("I77").Value = 38
$ws.Range("I78").Value = 96
$ws.Range("I79").Value = 191
$ws.Range("I83").Value = 140
$ws.Range("I84").Value = 51
$ws.Range("I85").Value = 347
$ws.Range("I89").Value = 75
$ws.Range("I96").Value = 96
$ws.Range("I98").Value = 51
$ws.Range("I99").Value = 134
$ws.Range("G101").Value = 24646
$ws.Range("I101").Value = 7390

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 51
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 84
$ws.Range("I4").Value = 20
$ws.Range("I7").Value = 349

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 35
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 167

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 144
$ws.Range("I3").Value = 159
$ws.Range("I7").Value = 477

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I3").Value = 134
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 347

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 89
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 19

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 53
$ws.Range("I3").Value = 55
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 16
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("G4").Value = 26
$ws.Range("G7").Value = 620

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("I2").Value = 2
$ws.Range("I7").Value = 8

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 144
$ws.Range("I5").Value = 16
$ws.Range("I6").Value = 156
$ws.Range("I7").Value = 469

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 15
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 33
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 69

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 25

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 14

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 25
